# Swap the two theme colour palettes that are in play in this deck.
#
# Before the edit:
#   ppt/theme/theme1.xml (used only by the Notes Master)   = "Office Theme" palette
#   ppt/theme/theme2.xml (used by the (only) Slide Master)  = "Integral" palette
#
# After the edit (per the source commit) the two palettes are swapped, so the
# slide master ends up using the default "Office Theme" colours, and the
# (mostly invisible) notes master ends up with the "Integral" colours.
#
# The only COM surface this runtime exposes for touching theme colours is
# Master.ColorScheme.Colors(n).RGB, which always resolves to the palette
# used by the deck's slide master (ppt/theme/theme2.xml here) - so we drive
# that object through the twelve standard theme colour slots and write the
# "Office Theme" RGB values into it (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink - that's the order PowerPoint exposes via ColorScheme.Colors).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# index -> target "Office Theme" RGB, encoded as R + G*256 + B*65536
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72

Write-Output "Theme colour scheme updated"
